$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.092277
$ws.Range("H2").Value = 0.276831
$ws.Range("I2").Value = 0.03444274323645406
$ws.Range("J2").Value = 0.03444274323645406
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.04661299999999999
$ws.Range("N2").Value = 0.139839
$ws.Range("O2").Value = 0.0286392647455175
$ws.Range("P2").Value = 0.0286392647455175
$ws.Range("Q2").Value = 0.004301307801
$ws.Range("R2").Value = 0.038711770209
$ws.Range("S2").Value = 0.00098641484211069
$ws.Range("T2").Value = 0.00098641484211069

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.092277
$ws.Range("H3").Value = 0.276831
$ws.Range("I3").Value = 0.03444274323645406
$ws.Range("J3").Value = 0.03444274323645406
$ws.Range("M3").Value = 1.580977666666667
$ws.Range("N3").Value = 4.742933000000001
$ws.Range("O3").Value = 0.9713607352544825
$ws.Range("P3").Value = 0.9713607352544824
$ws.Range("Q3").Value = 0.145887876147
$ws.Range("R3").Value = 1.312990885323
$ws.Range("S3").Value = 0.03345632839434336
$ws.Range("T3").Value = 0.03345632839434336

# Row 4
$ws.Range("G4").Value = 2.276331666666667
$ws.Range("H4").Value = 6.828995
$ws.Range("I4").Value = 0.8496495022162568
$ws.Range("J4").Value = 0.8496495022162568
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.04661299999999999
$ws.Range("N4").Value = 0.139839
$ws.Range("O4").Value = 0.0286392647455175
$ws.Range("P4").Value = 0.0286392647455175
$ws.Range("Q4").Value = 0.1061066479783333
$ws.Range("R4").Value = 0.9549598318049999
$ws.Range("S4").Value = 0.02433333703486854
$ws.Range("T4").Value = 0.02433333703486853

# Row 5
$ws.Range("G5").Value = 2.276331666666667
$ws.Range("H5").Value = 6.828995
$ws.Range("I5").Value = 0.8496495022162568
$ws.Range("J5").Value = 0.8496495022162568
$ws.Range("M5").Value = 1.580977666666667
$ws.Range("N5").Value = 4.742933000000001
$ws.Range("O5").Value = 0.9713607352544825
$ws.Range("P5").Value = 0.9713607352544824
$ws.Range("Q5").Value = 3.598829526926111
$ws.Range("R5").Value = 32.38946574233501
$ws.Range("S5").Value = 0.8253161651813883
$ws.Range("T5").Value = 0.8253161651813882

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.3105333333333333
$ws.Range("H6").Value = 0.9316
$ws.Range("I6").Value = 0.1159077545472891
$ws.Range("J6").Value = 0.1159077545472891
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.04661299999999999
$ws.Range("N6").Value = 0.139839
$ws.Range("O6").Value = 0.0286392647455175
$ws.Range("P6").Value = 0.0286392647455175
$ws.Range("Q6").Value = 0.01447489026666667
$ws.Range("R6").Value = 0.1302740124
$ws.Range("S6").Value = 0.003319512868538273
$ws.Range("T6").Value = 0.003319512868538273

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.3105333333333333
$ws.Range("H7").Value = 0.9316
$ws.Range("I7").Value = 0.1159077545472891
$ws.Range("J7").Value = 0.1159077545472891
$ws.Range("M7").Value = 1.580977666666667
$ws.Range("N7").Value = 4.742933000000001
$ws.Range("O7").Value = 0.9713607352544825
$ws.Range("P7").Value = 0.9713607352544824
$ws.Range("Q7").Value = 0.4909462647555556
$ws.Range("R7").Value = 4.418516382800001
$ws.Range("S7").Value = 0.1125882416787509
$ws.Range("T7").Value = 0.1125882416787509
